$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Primer Trimestre"
$ws.Range("B1").Value = 4

$ws.Range("B1").Select()
